$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" on Overview (G2) and "Correspond Handoff Datetime"
# on de-de (H2) share the same underlying text value, so both must be updated to the
# same new value together.
$wsOverview.Range("G2").Value = "2016-11-29 03:21:23"
$wsDeDe.Range("H2").Value = "2016-11-29 03:21:23"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsZhCn.Range("H2").Value = "2016-11-29 03:21:09"
$wsZhCn.Range("K2").Value = "2016-11-29 03:22:00"

# de-de sheet: "Correspond Handback DateTime" (K2)
$wsDeDe.Range("K2").Value = "2016-11-29 03:22:18"
